# Apply "commits ob and to" changes to the Official Business / Travel Order form.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- First copy of the form (rows 10-21) ---
# No. : 2020-007 -> 2020-006
$ws.Range("J12").Value = "2020-006"
# Date: (header) : May 04, 2020 -> April 30, 2020
$ws.Range("J13").Value = "April 30, 2020"
# purpose (s). value : RO -> Test for system
$ws.Range("C17").Value = "Test for system"
# Place to be visited: value : RO -> Regional Office
$ws.Range("D19").Value = "Regional Office"
# Date: (second occurrence, travel date) : May 01, 2020 -> April 30, 2020
$ws.Range("D21").Value = "April 30, 2020"
# Time of Return: value : 18:00:00 -> 17:00:00
$ws.Range("K21").Value = "17:00:00"

# --- Second copy of the form (rows 49-58, mirrors the first) ---
# No. : 2020-007 -> 2020-006
$ws.Range("J49").Value = "2020-006"
# Date: (header) : May 04, 2020 -> April 30, 2020
$ws.Range("J50").Value = "April 30, 2020"
# purpose (s). value : RO -> Test for system
$ws.Range("C54").Value = "Test for system"
# Place to be visited: value : RO -> Regional Office
$ws.Range("D56").Value = "Regional Office"
# Date: (second occurrence, travel date) : May 01, 2020 -> April 30, 2020
$ws.Range("D58").Value = "April 30, 2020"
# Time of Return: value : 18:00:00 -> 17:00:00
$ws.Range("K58").Value = "17:00:00"
